$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user retyped the "h" (height) column for every sprite row to the
# new uniform value of 121. Cells E8 and E9 previously held formulas
# (=E6 / =E7); typing a literal over them replaces the formula with the
# literal value, same as the other rows.
$ws.Range("E2:E87").Value = 121

# The first block's "y" column (C4:C5) had been hand-entered to mirror the
# old height (100) for the second visual row of that sprite grid; it was
# updated to match the new height as well.
$ws.Range("C4:C5").Value = 121

# A small manual correction to the "w" (width) of "caminar DER 1" (D22),
# from 112 to 111.
$ws.Range("D22").Value = 111

# Restore the view state: scrolled down with C6 selected.
$ws.Range("C6").Select()
